$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.224535333333333
$ws.Range("H2").Value = 12.673606
$ws.Range("I2").Value = 0.7043225486309714
$ws.Range("J2").Value = 0.7043225486309715
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.847498666666667
$ws.Range("N2").Value = 14.542496
$ws.Range("O2").Value = 0.03400671694637637
$ws.Range("P2").Value = 0.03400671694637637
$ws.Range("Q2").Value = 20.47842939561955
$ws.Range("R2").Value = 184.305864560576
$ws.Range("S2").Value = 0.02395169755024385
$ws.Range("T2").Value = 0.02395169755024385
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.224535333333333
$ws.Range("H3").Value = 12.673606
$ws.Range("I3").Value = 0.7043225486309714
$ws.Range("J3").Value = 0.7043225486309715
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.627093333333334
$ws.Range("N3").Value = 4.88128
$ws.Range("O3").Value = 0.01141456784970118
$ws.Range("P3").Value = 0.01141456784970118
$ws.Range("Q3").Value = 6.873713277297778
$ws.Range("R3").Value = 61.86341949568
$ws.Range("S3").Value = 0.008039537519422683
$ws.Range("T3").Value = 0.008039537519422685
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.224535333333333
$ws.Range("H4").Value = 12.673606
$ws.Range("I4").Value = 0.7043225486309714
$ws.Range("J4").Value = 0.7043225486309715
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 136.0707373333333
$ws.Range("N4").Value = 408.212212
$ws.Range("O4").Value = 0.9545787152039225
$ws.Range("P4").Value = 0.9545787152039225
$ws.Range("Q4").Value = 574.8356376973858
$ws.Range("R4").Value = 5173.520739276472
$ws.Range("S4").Value = 0.6723313135613049
$ws.Range("T4").Value = 0.672331313561305
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.773477
$ws.Range("H5").Value = 5.320431
$ws.Range("I5").Value = 0.2956774513690286
$ws.Range("J5").Value = 0.2956774513690286
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.847498666666667
$ws.Range("N5").Value = 14.542496
$ws.Range("O5").Value = 0.03400671694637637
$ws.Range("P5").Value = 0.03400671694637637
$ws.Range("Q5").Value = 8.596927392864
$ws.Range("R5").Value = 77.372346535776
$ws.Range("S5").Value = 0.01005501939613252
$ws.Range("T5").Value = 0.01005501939613252
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.773477
$ws.Range("H6").Value = 5.320431
$ws.Range("I6").Value = 0.2956774513690286
$ws.Range("J6").Value = 0.2956774513690286
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.627093333333334
$ws.Range("N6").Value = 4.88128
$ws.Range("O6").Value = 0.01141456784970118
$ws.Range("P6").Value = 0.01141456784970118
$ws.Range("Q6").Value = 2.88561260352
$ws.Range("R6").Value = 25.97051343168
$ws.Range("S6").Value = 0.003375030330278498
$ws.Range("T6").Value = 0.003375030330278498
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.773477
$ws.Range("H7").Value = 5.320431
$ws.Range("I7").Value = 0.2956774513690286
$ws.Range("J7").Value = 0.2956774513690286
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 136.0707373333333
$ws.Range("N7").Value = 408.212212
$ws.Range("O7").Value = 0.9545787152039225
$ws.Range("P7").Value = 0.9545787152039225
$ws.Range("Q7").Value = 241.318323033708
$ws.Range("R7").Value = 2171.864907303372
$ws.Range("S7").Value = 0.2822474016426175
$ws.Range("T7").Value = 0.2822474016426175
